# Notification OR properties: reorder/fix the VT200_0009 & VT200_0013 JS
# test-step blocks (drop a redundant wait(2) before Lock_UnlockScreen
# (unlock)) and move them later in the step sequence; also swap the
# launch_App_Device(...RhodesActivity) calls for press_Key(Back) in the
# VT200_0012 / VT200_0014 steps, tweak the selected cell & a couple of
# column/row display settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Row 7 : VT200-0009 --------------------------------------------------
$ws.Range("G7").Value = @"
wait(3);
validate1;
link_Click(Application_test_link);
validate2;
SelectTestToRun(VT200_0009_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
Lock_UnlockScreen(lock);
Lock_UnlockScreen(unlock);
wait(2);
link_Click(stopCallback_name_xpath);
wait(2);
validate4;
"@

# --- Row 8 : VT200-0011 ---------------------------------------------------
$ws.Range("E8").Value = "Check for Application event callback with pressing power button after minimizing application"

$ws.Range("G8").Value = @"
wait(3);
validate1;
link_Click(Application_test_link);
validate2;
SelectTestToRun(VT200_0011_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
press_Key(home);
validate4;
Lock_UnlockScreen(lock);
Lock_UnlockScreen(unlock);
launch_App_Device(com.rhomobile.compliancetest_js/com.rhomobile.rhodes.RhodesActivity);
link_Click(stopCallback_name_xpath);
validate5;
"@

$ws.Range("H8").Value = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Application JS Test
};
validate3
{
validate_Text_Exists=VT200-0011
};
validate4
{
validate_AppMinimized=AppsScreen
};
validate5
{
validate_Result=UIDestroyed
validate_Result=Deactivated
validate_Result=Activated
validate_doesNotContain=ScreenOff
validate_doesNotContain=ScreenOn
};
"@

# --- Row 9 : VT200-0012 ---------------------------------------------------
$ws.Range("E9").Value = "Check for Application event callback with pressing power button after sending the app to background"

$ws.Range("G9").Value = @"
wait(3);
validate1;
link_Click(Application_test_link);
validate2;
SelectTestToRun(VT200_0012_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
launch_App_Device(com.android.browser/com.android.browser.BrowserActivity);
Lock_UnlockScreen(lock);
Lock_UnlockScreen(unlock);
press_Key(Back);
validate4;
link_Click(stopCallback_name_xpath);
wait(2);
validate5;
"@

$ws.Range("H9").Value = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Application JS Test
};
validate3
{
validate_Text_Exists=VT200-0012
};
validate4
{
validate_PageTitle=Application JS Test
};
validate5
{
validate_Result=UIDestroyed
validate_Result=Deactivated
validate_Result=Activated
validate_doesNotContain=ScreenOff
validate_doesNotContain=ScreenOn
};
"@

# --- Row 10 : VT200-0013 ---------------------------------------------------
$ws.Range("G10").Value = @"
wait(3);
validate1;
link_Click(Application_test_link);
validate2;
SelectTestToRun(VT200_0013_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(8);
Lock_UnlockScreen(lock);
Lock_UnlockScreen(unlock);
link_Click(stopCallback_name_xpath);
wait(2);
validate6;
"@

$ws.Range("H10").Value = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Application JS Test
};
validate3
{
validate_Text_Exists=VT200-0013
};
validate4
{
validate_AppMinimized=AppsScreen
};
validate5
{
validate_PageTitle=Application JS Test
};
validate6
{
validate_Result=UIDestroyed
validate_Result=Deactivated
validate_Result=ScreenOff
validate_Result=Activated
validate_Result=ScreenOn
};
"@

# --- Row 11 : VT200-0014 ---------------------------------------------------
$ws.Range("G11").Value = @"
wait(3);
validate1;
link_Click(Application_test_link);
validate2;
SelectTestToRun(VT200_0014_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
launch_App_Device(com.android.browser/com.android.browser.BrowserActivity);
press_Key(Back);
Lock_UnlockScreen(lock);
Lock_UnlockScreen(unlock);
wait(2);
link_Click(stopCallback_name_xpath);
wait(2);
validate4;
"@

$ws.Range("H11").Value = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Application JS Test
};
validate3
{
validate_Text_Exists=VT200-0014
};
validate4
{
validate_Result=UIDestroyed
validate_Result=Deactivated
validate_Result=ScreenOff
validate_Result=Activated
validate_Result=ScreenOn
};
"@

$ws.Rows(11).RowHeight = 255.75

# --- Row 12 : VT200-0016 ---------------------------------------------------
$ws.Range("G12").Value = @"
wait(5);
validate1;
link_Click(Application_test_link);
wait(5);
validate2;
SelectTestToRun(VT200_0016_string);
ClickRunTest(runtest_top_xpath);
wait(5);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
press_Menu(menu);
wait(2);
ClickUITextView(Load_Page);
wait(2);
validate4;

"@

$ws.Range("H12").Value = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Application JS Test
};
validate3
{
validate_Text_Exists=VT200-0016
};
validate4
{
validate_Page=Loading
};

"@

# --- sheet view selection ---------------------------------------------
$ws.Range("H2").Select()

# --- column display flags (TestCases!E, Sheet2 A/C/D/F) ----------------
$ws.Columns("E").ColumnWidth = 19

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns("A").ColumnWidth = 17.85546875
$ws2.Columns("C").ColumnWidth = 23.140625
$ws2.Columns("D").ColumnWidth = 34.5703125
$ws2.Columns("F").ColumnWidth = 18.28515625
